# "Ultimo excel de socios y empresas"
#
# 1) Registro de usuarios (sheet1): drop the "Username" column, retype the
#    roster with new people/rows, re-point the mail hyperlinks.
# 2) Add a new "Registro de empresas" sheet (sheet2) with a company-registry
#    header row, and make it the active tab.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# Sheet1 ("Registro de usuarios")
# ---------------------------------------------------------------------

# Drop the old hyperlink (D2 -> alfonsi@gmail.com) entirely; the column it
# lived in is going away.
$ws1.Range("D2").Hyperlinks.Delete()

# Headers: A stays "Nomre y apellidos"; B..E shift in from the old C..F,
# losing the "Username" header (old B) altogether.
$ws1.Range("B1").Value = "Empresa"
$ws1.Range("C1").Value = "Mail"
$ws1.Range("D1").Value = "Username Generado"
$ws1.Range("E1").Value = "PW Generado"
$ws1.Range("F1").Clear()

# Wipe the old data row so stale cells (old C2/D2) don't linger once the
# new, narrower rows are written.
$ws1.Range("A2:F2").Clear()

# Row 2 - Goberto Calleja Calleja / Secpho
$ws1.Range("A2").Value = "Goberto Calleja Calleja"
$ws1.Range("B2").Value = "Secpho"
$ws1.Range("C2").Value = "user1@gmail.com"
$ws1.Range("D2").Value = "user1"
$ws1.Range("E2").Value = 12345

# Row 3 - Victor Gonzales / Secpho
$ws1.Range("A3").Value = "Victor Gonzales"
$ws1.Range("B3").Value = "Secpho"
$ws1.Range("C3").Value = "VicG@gmail.com"

# Row 4 - Daniel Carvajal / Inventado
$ws1.Range("A4").Value = "Daniel Carvajal"
$ws1.Range("B4").Value = "Inventado"
$ws1.Range("C4").Value = "danic@gmail.com"

# Mail column hyperlinks (now column C)
$ws1.Hyperlinks.Add($ws1.Range("C2"), "mailto:user1@gmail.com")
$ws1.Range("C2").Style = "Hyperlink"
$ws1.Hyperlinks.Add($ws1.Range("C3"), "mailto:VicG@gmail.com")
$ws1.Range("C3").Style = "Hyperlink"
$ws1.Hyperlinks.Add($ws1.Range("C4"), "mailto:danic@gmail.com")
$ws1.Range("C4").Style = "Hyperlink"

# Column C used to be unstyled/narrow; it now takes over the old Mail
# column's width.
$ws1.Range("C1").ColumnWidth = 19.5

# Leave the cursor parked below the data, like the source file does.
$ws1.Range("A7").Select()

# ---------------------------------------------------------------------
# Sheet2 ("Registro de empresas") - new
# ---------------------------------------------------------------------

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Registro de empresas"

# Reuse sheet1's bold/filled header formatting instead of rebuilding it
# from scratch (keeps the shared style table free of near-duplicates).
$ws1.Range("A1:E1").Copy()
$ws2.Range("A1:H1").PasteSpecial(-4122)

$ws2.Range("A1").Value = "Name"
$ws2.Range("B1").Value = "Descripción"
$ws2.Range("C1").Value = "Url"
$ws2.Range("D1").Value = "Tipo De Entidad"
$ws2.Range("E1").Value = "Actividad"
$ws2.Range("F1").Value = "Dirección"
$ws2.Range("G1").Value = "Contacto"
$ws2.Range("H1").Value = "Logo Url"

$ws2.Range("B1").ColumnWidth = 12.1
$ws2.Range("D1").ColumnWidth = 16

# Select the header row and make this the active/visible sheet, matching
# the saved view state.
$ws2.Range("A1:H1").Select()
$ws2.Activate()
